$wb = $excel.ActiveWorkbook

# --- Sheet: Neg_Change ---
$ws1 = $wb.Worksheets.Item("Neg_Change")

# Remove rows 5-10 (old extra tickers); new data only has rows 2-4
$ws1.Range("A5:I10").EntireRow.Delete()

# Row 2: ONGC
$ws1.Range("A2").Value = "ONGC"
$ws1.Range("B2").Value = 234.95
$ws1.Range("C2").Value = 235.6
$ws1.Range("D2").Value = 233.31
$ws1.Range("E2").Value = 234.5
$ws1.Range("F2").Value = 6236729
$ws1.Range("G2").Value = 12621499
$ws1.Range("H2").Value = -0.5058646362052558
$ws1.Range("I2").Value = "ONGC"

# Row 3: HINDZINC
$ws1.Range("A3").Value = "HINDZINC"
$ws1.Range("B3").Value = 615
$ws1.Range("C3").Value = 630
$ws1.Range("D3").Value = 613.05
$ws1.Range("E3").Value = 625.85
$ws1.Range("F3").Value = 11718688
$ws1.Range("G3").Value = 27992946
$ws1.Range("H3").Value = -0.5813699637044275
$ws1.Range("I3").Value = "HINDZINC"

# Row 4: GRANULES
$ws1.Range("A4").Value = "GRANULES"
$ws1.Range("B4").Value = 601
$ws1.Range("C4").Value = 604.95
$ws1.Range("D4").Value = 594.8
$ws1.Range("E4").Value = 598
$ws1.Range("F4").Value = 427337
$ws1.Range("G4").Value = 899319
$ws1.Range("H4").Value = -0.5248215594244089
$ws1.Range("I4").Value = "GRANULES"

# --- Sheet: Pos_Change ---
$ws2 = $wb.Worksheets.Item("Pos_Change")

# Row 2: HINDALCO
$ws2.Range("A2").Value = "HINDALCO"
$ws2.Range("B2").Value = 868
$ws2.Range("C2").Value = 888.3
$ws2.Range("D2").Value = 863.85
$ws2.Range("E2").Value = 883.3
$ws2.Range("F2").Value = 15280248
$ws2.Range("G2").Value = 9836341
$ws2.Range("H2").Value = 0.5534483808562554
$ws2.Range("I2").Value = "HINDALCO"

# Row 3: LT
$ws2.Range("A3").Value = "LT"
$ws2.Range("B3").Value = 4039
$ws2.Range("C3").Value = 4063.9
$ws2.Range("D3").Value = 4006.2
$ws2.Range("E3").Value = 4056
$ws2.Range("F3").Value = 1362208
$ws2.Range("G3").Value = 952143
$ws2.Range("H3").Value = 0.4306758543622124
$ws2.Range("I3").Value = "LT"

# Row 4: SUNPHARMA
$ws2.Range("A4").Value = "SUNPHARMA"
$ws2.Range("B4").Value = 1723.1
$ws2.Range("C4").Value = 1723.1
$ws2.Range("D4").Value = 1707
$ws2.Range("E4").Value = 1720
$ws2.Range("F4").Value = 2639497
$ws2.Range("G4").Value = 1800816
$ws2.Range("H4").Value = 0.4657227612371281
$ws2.Range("I4").Value = "SUNPHARMA"

# Row 5: HINDUNILVR
$ws2.Range("A5").Value = "HINDUNILVR"
$ws2.Range("B5").Value = 2293.3
$ws2.Range("C5").Value = 2298
$ws2.Range("D5").Value = 2283
$ws2.Range("E5").Value = 2290
$ws2.Range("F5").Value = 1556717
$ws2.Range("G5").Value = 1086695
$ws2.Range("H5").Value = 0.4325243053478667
$ws2.Range("I5").Value = "HINDUNILVR"

# Row 6: JIOFIN
$ws2.Range("A6").Value = "JIOFIN"
$ws2.Range("B6").Value = 293
$ws2.Range("C6").Value = 293.85
$ws2.Range("D6").Value = 291.4
$ws2.Range("E6").Value = 293.2
$ws2.Range("F6").Value = 6850674
$ws2.Range("G6").Value = 4713140
$ws2.Range("H6").Value = 0.4535265237187947
$ws2.Range("I6").Value = "JIOFIN"

# Row 7: RELIANCE
$ws2.Range("A7").Value = "RELIANCE"
$ws2.Range("B7").Value = 1547
$ws2.Range("C7").Value = 1553.6
$ws2.Range("D7").Value = 1537.8
$ws2.Range("E7").Value = 1541
$ws2.Range("F7").Value = 8815884
$ws2.Range("G7").Value = 5972105
$ws2.Range("H7").Value = 0.4761769928693484
$ws2.Range("I7").Value = "RELIANCE"

# Row 8: TRENT
$ws2.Range("A8").Value = "TRENT"
$ws2.Range("B8").Value = 4226
$ws2.Range("C8").Value = 4257.7
$ws2.Range("D8").Value = 4187.8
$ws2.Range("E8").Value = 4210
$ws2.Range("F8").Value = 670540
$ws2.Range("G8").Value = 423139
$ws2.Range("H8").Value = 0.5846802114671538
$ws2.Range("I8").Value = "TRENT"

# Row 9: POWERGRID
$ws2.Range("A9").Value = "POWERGRID"
$ws2.Range("B9").Value = 260
$ws2.Range("C9").Value = 261.45
$ws2.Range("D9").Value = 258.95
$ws2.Range("E9").Value = 259.45
$ws2.Range("F9").Value = 10929200
$ws2.Range("G9").Value = 7270521
$ws2.Range("H9").Value = 0.5032210208869489
$ws2.Range("I9").Value = "POWERGRID"

# Row 10: BANKBARODA
$ws2.Range("A10").Value = "BANKBARODA"
$ws2.Range("B10").Value = 286.25
$ws2.Range("C10").Value = 293.8
$ws2.Range("D10").Value = 285.5
$ws2.Range("E10").Value = 292.9
$ws2.Range("F10").Value = 11505161
$ws2.Range("G10").Value = 7261978
$ws2.Range("H10").Value = 0.5843012743910819
$ws2.Range("I10").Value = "BANKBARODA"

# Row 11: GAIL
$ws2.Range("A11").Value = "GAIL"
$ws2.Range("B11").Value = 170.35
$ws2.Range("C11").Value = 171
$ws2.Range("D11").Value = 169.76
$ws2.Range("E11").Value = 170.9
$ws2.Range("F11").Value = 4526325
$ws2.Range("G11").Value = 2857519
$ws2.Range("H11").Value = 0.5840052157133513
$ws2.Range("I11").Value = "GAIL"

# Row 12: PIDILITIND
$ws2.Range("A12").Value = "PIDILITIND"
$ws2.Range("B12").Value = 1452.8
$ws2.Range("C12").Value = 1457.1
$ws2.Range("D12").Value = 1437.7
$ws2.Range("E12").Value = 1455
$ws2.Range("F12").Value = 832687
$ws2.Range("G12").Value = 587342
$ws2.Range("H12").Value = 0.4177208508841527
$ws2.Range("I12").Value = "PIDILITIND"

# Row 13: ZYDUSLIFE
$ws2.Range("A13").Value = "ZYDUSLIFE"
$ws2.Range("B13").Value = 904
$ws2.Range("C13").Value = 906
$ws2.Range("D13").Value = 896.75
$ws2.Range("E13").Value = 902.4
$ws2.Range("F13").Value = 721717
$ws2.Range("G13").Value = 487836
$ws2.Range("H13").Value = 0.4794254626554826
$ws2.Range("I13").Value = "ZYDUSLIFE"

# Row 14: CGPOWER
$ws2.Range("A14").Value = "CGPOWER"
$ws2.Range("B14").Value = 646.9
$ws2.Range("C14").Value = 648
$ws2.Range("D14").Value = 637.05
$ws2.Range("E14").Value = 643.1
$ws2.Range("F14").Value = 1954985
$ws2.Range("G14").Value = 1234659
$ws2.Range("H14").Value = 0.5834210093637191
$ws2.Range("I14").Value = "CGPOWER"

# Row 15: OIL
$ws2.Range("A15").Value = "OIL"
$ws2.Range("B15").Value = 408.2
$ws2.Range("C15").Value = 413
$ws2.Range("D15").Value = 404.8
$ws2.Range("E15").Value = 412.1
$ws2.Range("F15").Value = 1277438
$ws2.Range("G15").Value = 893809
$ws2.Range("H15").Value = 0.4292069110962186
$ws2.Range("I15").Value = "OIL"

# Row 16: BSE
$ws2.Range("A16").Value = "BSE"
$ws2.Range("B16").Value = 2628.3
$ws2.Range("C16").Value = 2629.7
$ws2.Range("D16").Value = 2570.2
$ws2.Range("E16").Value = 2595
$ws2.Range("F16").Value = 3829439
$ws2.Range("G16").Value = 2577976
$ws2.Range("H16").Value = 0.4854440072366849
$ws2.Range("I16").Value = "BSE"

# Row 17: TATATECH
$ws2.Range("A17").Value = "TATATECH"
$ws2.Range("B17").Value = 650.1
$ws2.Range("C17").Value = 651.2
$ws2.Range("D17").Value = 638
$ws2.Range("E17").Value = 640
$ws2.Range("F17").Value = 614244
$ws2.Range("G17").Value = 408887
$ws2.Range("H17").Value = 0.5022341135815029
$ws2.Range("I17").Value = "TATATECH"

# Row 18: TATAELXSI
$ws2.Range("A18").Value = "TATAELXSI"
$ws2.Range("B18").Value = 5318
$ws2.Range("C18").Value = 5342
$ws2.Range("D18").Value = 5175.5
$ws2.Range("E18").Value = 5210.5
$ws2.Range("F18").Value = 141898
$ws2.Range("G18").Value = 96723
$ws2.Range("H18").Value = 0.4670554056429184
$ws2.Range("I18").Value = "TATAELXSI"

# Row 19: BANDHANBNK
$ws2.Range("A19").Value = "BANDHANBNK"
$ws2.Range("B19").Value = 146.15
$ws2.Range("C19").Value = 146.6
$ws2.Range("D19").Value = 144.73
$ws2.Range("E19").Value = 145.6
$ws2.Range("F19").Value = 5408553
$ws2.Range("G19").Value = 3804621
$ws2.Range("H19").Value = 0.4215747113838672
$ws2.Range("I19").Value = "BANDHANBNK"

# Row 20: CAMS
$ws2.Range("A20").Value = "CAMS"
$ws2.Range("B20").Value = 738.7
$ws2.Range("C20").Value = 740.1
$ws2.Range("D20").Value = 729.7
$ws2.Range("E20").Value = 735.4
$ws2.Range("F20").Value = 792833
$ws2.Range("G20").Value = 534717
$ws2.Range("H20").Value = 0.4827151558674962
$ws2.Range("I20").Value = "CAMS"

# Row 21: DELHIVERY
$ws2.Range("A21").Value = "DELHIVERY"
$ws2.Range("B21").Value = 403.1
$ws2.Range("C21").Value = 407.3
$ws2.Range("D21").Value = 399.95
$ws2.Range("E21").Value = 401.5
$ws2.Range("F21").Value = 1244456
$ws2.Range("G21").Value = 872556
$ws2.Range("H21").Value = 0.4262190621576151
$ws2.Range("I21").Value = "DELHIVERY"
